# Updated cryptos list on Mon Jan 29 17:40:41 UTC 2024 with GitHub Actions
# Refresh the price / 1h-volume columns (and a couple of reordered rows) on the
# crypto ranking sheet to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.137.14'
$ws.Range('E2').Value = '  +2.38%  '
$ws.Range('D3').Value = '2.306.44'
$ws.Range('E3').Value = '  +1.89%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''309.75'
$ws.Range('E5').Value = '  +1.54%  '
$ws.Range('D6').Value = '''101.00'
$ws.Range('E6').Value = '  +5.33%  '
$ws.Range('D7').Value = '''0.537'
$ws.Range('E7').Value = '  +1.46%  '
$ws.Range('D9').Value = '''0.510'
$ws.Range('E9').Value = '  +4.23%  '
$ws.Range('D10').Value = '''36.13'
$ws.Range('E10').Value = '  +2.80%  '
$ws.Range('D11').Value = '''0.0821'
$ws.Range('E11').Value = '  +3.61%  '
$ws.Range('E12').Value = '  +0.88%  '
$ws.Range('D13').Value = '''6.98'
$ws.Range('E13').Value = '  +5.38%  '
$ws.Range('D14').Value = '2.661.03'
$ws.Range('E14').Value = '  +1.60%  '
$ws.Range('D15').Value = '''14.92'
$ws.Range('E15').Value = '  +3.98%  '
$ws.Range('D16').Value = '2.295.58'
$ws.Range('E16').Value = '  +0.79%  '
$ws.Range('D17').Value = '''0.806'
$ws.Range('E17').Value = '  +1.84%  '
$ws.Range('D18').Value = '43.091.19'
$ws.Range('E18').Value = '  +2.43%  '
$ws.Range('D19').Value = '''12.60'
$ws.Range('E19').Value = '  +1.27%  '
$ws.Range('D20').Value = '0.0₃0919'
$ws.Range('E20').Value = '  +1.44%  '
$ws.Range('E21').Value = '  +1.81%  '
$ws.Range('D22').Value = '''68.28'
$ws.Range('E22').Value = '  +0.85%  '
$ws.Range('D23').Value = '''240.35'
$ws.Range('E23').Value = '  +1.10%  '
$ws.Range('E24').Value = '  +5.10%  '
$ws.Range('E25').Value = '  +2.17%  '
$ws.Range('D26').Value = '''0.998'
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').Value = '''24.16'
$ws.Range('E27').Value = '  +1.97%  '
$ws.Range('D28').Value = '''38.91'
$ws.Range('E28').Value = '  +5.57%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').Value = '''9.68'
$ws.Range('E29').Value = '  +1.86%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '''2.13'
$ws.Range('E30').Value = '  +0.82%  '
$ws.Range('D31').Value = '''165.55'
$ws.Range('E31').Value = '  +3.87%  '
$ws.Range('E32').Value = '  +1.72%  '
$ws.Range('D33').Value = '''1.00'
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').Value = '''3.18'
$ws.Range('E34').Value = '  -0.29%  '
$ws.Range('D35').Value = '''17.85'
$ws.Range('E35').Value = '  +4.62%  '
$ws.Range('D36').Value = '''0.0742'
$ws.Range('E36').Value = '  +0.62%  '
$ws.Range('E37').Value = '  +0.58%  '
$ws.Range('E38').Value = '  +0.64%  '
$ws.Range('E39').Value = '  +1.64%  '
$ws.Range('E40').Value = '  +1.71%  '
$ws.Range('D41').Value = '''4.20'
$ws.Range('E41').Value = '  +3.91%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '''19.72'
$ws.Range('E42').Value = '  +4.66%  '
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').Value = '''2.32'
$ws.Range('E43').Value = '  -5.36%  '
$ws.Range('D44').Value = '''0.0291'
$ws.Range('E44').Value = '  +2.47%  '
$ws.Range('D45').Value = '1.969.60'
$ws.Range('E45').Value = '  -0.88%  '
$ws.Range('E46').Value = '  +4.02%  '
$ws.Range('D47').Value = '''9.84'
$ws.Range('E47').Value = '  -0.81%  '
$ws.Range('D48').Value = '''3.01'
$ws.Range('E48').Value = '  +20.61%  '
$ws.Range('D49').Value = '''55.19'
$ws.Range('E49').Value = '  +3.92%  '
$ws.Range('D50').Value = '''1.55'
$ws.Range('E50').Value = '  +3.01%  '
$ws.Range('D51').Value = '2.528.11'
$ws.Range('E51').Value = '  +1.55%  '
